$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "metformin-pioglitazone"
$ws.Range("A28").Value = "pioglitazone-metformin"
$ws.Range("A29").Value = "duetact"

$fontRange = $ws.Range("A27:A28")
$fontRange.Font.Name = "Lucida Console"
$fontRange.Font.Size = 7
$fontRange.Font.Color = 0
$fontRange.VerticalAlignment = -4108
